# Core 64 Teensy LC Pin Usage Mapping - layout updates (V0.1, 75% complete)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column P (pin-specific notes, right side of table) ---

# Row 6 / P6: clarify that the 3.3V regulated VOUT also feeds "low power stuff"
$ws.Range("P6").Value = "3.3V Regulated VOUT to AREF and low power stuff"

# Rows 11 & 12 / P11, P12: shorten the OLED/hall-sensor note, dropping the open question
$ws.Range("P11").Value = "OLED Debug screen, hall sensors."
$ws.Range("P12").Value = "OLED Debug screen, hall sensors."

# --- Rows 21-23 (column A / I / P): assign previously unused / "No Connection" pins ---

# Row 21 / A21: was "No Connection" -> now used for LOGIC Ground
$ws.Range("A21").Value = "LOGIC Ground"

# Row 22 / A22, I22, P22: was "No Connection" -> now documents a pull-high / VUSB power path
$ws.Range("A22").Value = "Pull high to avoid acceptable programming"
$ws.Range("I22").Value = "VUSB"
$ws.Range("P22").Value = "Power from USB to LiPo charger"

# Row 23 / A23: was "Available" -> now assigned to the Hall Switch
$ws.Range("A23").Value = "Hall Switch"

# --- View state: selection moved from A24 to A22, scrolled back to top of sheet ---
$ws.Range("A22").Select()
